$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-changed date for each logging
# notification row. Rows 2-14 all currently show serial date 45221
# (2023-10-22); the automatic update refreshes them to serial date 45224
# (2023-10-25), leaving everything else (style, other columns) untouched.
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
